# petty-cashBook-2021.xlsx — "Update 19-Apr-2021, midday update."
#
# Target sheet is "Sheet1" (the first / active tab, containing the daily
# petty-cash ledger rows 1-113).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 35: add a second 300,000 disbursement on top of the existing 60,000 ---
$ws.Range("D35").Formula = "=60000+300000"

# --- Row 36: add two more components (900,000 + 1,209,000) to the existing debit ---
$ws.Range("D36").Formula = "=5000000+9750000+900000+1209000"

# --- Row 37: add a 51,187,000 credit on top of the existing 9,750,000 ---
$ws.Range("C37").Formula = "=9750000+51187000"

# --- Row 39: new transaction — "PAPA - beli obat", 350,000 debit ---
$ws.Range("B39").Value = "PAPA - beli obat"
$ws.Range("D39").Value = 350000

# --- Row 40: new transaction — "SALES - cash/retail" credit ---
$ws.Range("B40").Value = "SALES - cash/retail"
$ws.Range("C40").Formula = "=59197025+6708975-51187000"

# --- Row 41: new transaction — "SELISIH - lebih" credit of 95,000 ---
$ws.Range("B41").Value = "SELISIH - lebih"
$ws.Range("C41").Value = 95000

# --- Row 42: new transaction — "SETOR KE BANK" debit of 58,000,000 ---
$ws.Range("B42").Value = "SETOR KE BANK"
$ws.Range("D42").Value = 58000000

# --- Row 43: new date entry, 19-Apr-2021 (serial 44303) ---
$ws.Range("A43").Value = 44303

# Move the view: the frozen header stays at row 2, but the visible window and
# the active selection both move down to around the newly-added rows.
$ws.Range("B43").Select()
